$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 79245
$ws.Range("A7").Value = 130961458
$ws.Range("B7").Value = 79245
$ws.Range("Q7").Value = 446059
$ws.Range("R7").Value = 6760088
$ws.Range("A8").Value = 130961962
$ws.Range("B8").Value = 79245
$ws.Range("Q8").Value = 446084
$ws.Range("R8").Value = 6759981
$ws.Range("B9").Value = 79245
$ws.Range("B11").Value = 79864
$ws.Range("A13").Value = 130961060
$ws.Range("B13").Value = 79245
$ws.Range("Q13").Value = 446138
$ws.Range("R13").Value = 6759967
$ws.Range("Z13").Value = "10:26"
$ws.Range("AB13").Value = "10:26"
$ws.Range("A14").Value = 130963950
$ws.Range("B14").Value = 79245
$ws.Range("Q14").Value = 445926
$ws.Range("R14").Value = 6760113
$ws.Range("Z14").Value = "14:08"
$ws.Range("AB14").Value = "14:08"
$ws.Range("B15").Value = 79245
$ws.Range("A16").Value = 130963873
$ws.Range("B16").Value = 79245
$ws.Range("Q16").Value = 445938
$ws.Range("R16").Value = 6760155
$ws.Range("B17").Value = 79245
$ws.Range("A18").Value = 130961956
$ws.Range("B18").Value = 79864
$ws.Range("E18").Value = 6453
$ws.Range("F18").Value = "Vedskivlav"
$ws.Range("G18").Value = "Hertelidea botryosa"
$ws.Range("H18").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q18").Value = 446084
$ws.Range("R18").Value = 6759981
$ws.Range("AC18").Value = "Miljöbilder"
$ws.Range("A19").Value = 130960843
$ws.Range("B19").Value = 79245
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("Q19").Value = 446247
$ws.Range("R19").Value = 6759903
$ws.Range("AC19").Value = ""
$ws.Range("B21").Value = 79245
$ws.Range("A22").Value = 130962722
$ws.Range("B22").Value = 79864
$ws.Range("E22").Value = 6453
$ws.Range("F22").Value = "Vedskivlav"
$ws.Range("G22").Value = "Hertelidea botryosa"
$ws.Range("H22").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q22").Value = 446008
$ws.Range("R22").Value = 6759948
$ws.Range("Z22").Value = "10:26"
$ws.Range("AB22").Value = "10:26"
$ws.Range("AC22").Value = ""
$ws.Range("A23").Value = 130962640
$ws.Range("B23").Value = 79864
$ws.Range("Q23").Value = 446038
$ws.Range("R23").Value = 6759945
$ws.Range("A24").Value = 130963976
$ws.Range("B24").Value = 79245
$ws.Range("E24").Value = 6425
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("Q24").Value = 445929
$ws.Range("R24").Value = 6760099
$ws.Range("Z24").Value = "14:08"
$ws.Range("AB24").Value = "14:08"
$ws.Range("AC24").Value = "Miljöbild"
$ws.Range("B25").Value = 79245
$ws.Range("A27").Value = 130962736
$ws.Range("B27").Value = 79835
$ws.Range("E27").Value = 229821
$ws.Range("F27").Value = "Vedflamlav"
$ws.Range("G27").Value = "Ramboldia elabens"
$ws.Range("H27").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M27").Value = ""
$ws.Range("Q27").Value = 446008
$ws.Range("R27").Value = 6759948
$ws.Range("Z27").Value = "10:26"
$ws.Range("AB27").Value = "10:26"
$ws.Range("A28").Value = 130961461
$ws.Range("B28").Value = 79245
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 446088
$ws.Range("R28").Value = 6760088
$ws.Range("A29").Value = 130961750
$ws.Range("B29").Value = 79245
$ws.Range("Q29").Value = 446098
$ws.Range("R29").Value = 6760061
$ws.Range("AC29").Value = "Rikligt i en radie av ca 50 meter"
$ws.Range("A30").Value = 130963807
$ws.Range("B30").Value = 57881
$ws.Range("E30").Value = 100049
$ws.Range("F30").Value = "Spillkråka"
$ws.Range("G30").Value = "Dryocopus martius"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "färska spår"
$ws.Range("Q30").Value = 445932
$ws.Range("R30").Value = 6760079
$ws.Range("Z30").Value = "14:08"
$ws.Range("AB30").Value = "14:08"
$ws.Range("AC30").Value = ""
$ws.Range("B31").Value = 79245
$ws.Range("B32").Value = 79245
$ws.Range("A33").Value = 130961219
$ws.Range("B33").Value = 79245
$ws.Range("Q33").Value = 446122
$ws.Range("R33").Value = 6760020
$ws.Range("A34").Value = 130962676
$ws.Range("B34").Value = 79245
$ws.Range("Q34").Value = 446038
$ws.Range("R34").Value = 6759945
